{"js": "// Replace the 25 division-problem answers in the practice table, in document\n// (reading) order. The table has 20 rows x 5 columns, but only every 4th row\n// (0, 4, 8, 12, 16) actually holds text -- the rows in between are blank\n// \"work space\" rows for the student. We walk the non-blank rows top to\n// bottom, left to right, and overwrite each cell's text with the new value.\nconst replacements = [\n  \"23\u00f73=7, 2\", \"76\u00f76=12, 4\", \"29\u00f79=3, 2\", \"64\u00f75=12, 4\", \"17\u00f76=2, 5\",\n  \"92\u00f76=15, 2\", \"22\u00f72=11, 0\", \"42\u00f74=10, 2\", \"49\u00f72=24, 1\", \"27\u00f77=3, 6\",\n  \"58\u00f78=7, 2\", \"24\u00f75=4, 4\", \"14\u00f77=2, 0\", \"83\u00f72=41, 1\", \"53\u00f77=7, 4\",\n  \"16\u00f72=8, 0\", \"83\u00f75=16, 3\", \"76\u00f77=10, 6\", \"93\u00f77=13, 2\", \"92\u00f74=23, 0\",\n  \"98\u00f76=16, 2\", \"37\u00f74=9, 1\", \"48\u00f79=5, 3\", \"10\u00f72=5, 0\", \"50\u00f76=8, 2\"\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount,values\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst values = table.values;\nconst colCount = 5;\n\nlet idx = 0;\nfor (let r = 0; r < rowCount && idx < replacements.length; r++) {\n  const rowValues = values[r] || [];\n  const hasText = rowValues.some((v) => v !== \"\");\n  if (!hasText) {\n    continue; // skip the blank \"work\" rows\n  }\n  for (let c = 0; c < colCount && idx < replacements.length; c++) {\n    table.getCell(r, c).value = replacements[idx];\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the practice table, in document\n# (reading) order. The table has 20 rows x 5 columns, but only every 4th row\n# (1, 5, 9, 13, 17 in 1-based COM indexing) actually holds text -- the rows\n# in between are blank \"work space\" rows for the student. We walk the table\n# row by row, column by column, skip blank rows, and overwrite each cell's\n# text with the new value.\n$replacements = @(\n  \"23\u00f73=7, 2\", \"76\u00f76=12, 4\", \"29\u00f79=3, 2\", \"64\u00f75=12, 4\", \"17\u00f76=2, 5\",\n  \"92\u00f76=15, 2\", \"22\u00f72=11, 0\", \"42\u00f74=10, 2\", \"49\u00f72=24, 1\", \"27\u00f77=3, 6\",\n  \"58\u00f78=7, 2\", \"24\u00f75=4, 4\", \"14\u00f77=2, 0\", \"83\u00f72=41, 1\", \"53\u00f77=7, 4\",\n  \"16\u00f72=8, 0\", \"83\u00f75=16, 3\", \"76\u00f77=10, 6\", \"93\u00f77=13, 2\", \"92\u00f74=23, 0\",\n  \"98\u00f76=16, 2\", \"37\u00f74=9, 1\", \"48\u00f79=5, 3\", \"10\u00f72=5, 0\", \"50\u00f76=8, 2\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$idx = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n  # Peek at the row's first cell to decide whether this row holds answers.\n  $firstCellText = $t.Cell($r, 1).Range.Text\n  $firstCellText = $firstCellText.TrimEnd([char]13, [char]7)\n  if ($firstCellText -eq \"\") {\n    continue\n  }\n  for ($c = 1; $c -le $t.Columns.Count; $c++) {\n    if ($idx -ge $replacements.Count) {\n      break\n    }\n    $t.Cell($r, $c).Range.Text = $replacements[$idx]\n    $idx++\n  }\n}\n"}
